# Updates cryptos list values per upstream data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '59.060.01'
$ws.Range("E2").Value = '  +1.29%  '
$ws.Range("D3").Value = '2.970.01'
$ws.Range("E3").Value = '  -0.84%  '
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").Value = "'" + '563.07'
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'" + '136.43'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.24%  '
$ws.Range("E7").Value = '  -0.07%  '
$ws.Range("E8").Value = '  -0.26%  '
$ws.Range("D9").Value = '2.965.28'
$ws.Range("E9").Value = '  -0.63%  '
$ws.Range("E10").Value = '  +0.85%  '
$ws.Range("D11").Value = "'" + '5.26'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +7.84%  '
$ws.Range("D12").Value = "'" + '0.447'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.65%  '
$ws.Range("D14").Value = "'" + '33.45'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.40%  '
$ws.Range("E15").Value = '  -0.52%  '
$ws.Range("D16").Value = '3.460.13'
$ws.Range("E16").Value = '  -0.80%  '
$ws.Range("E17").Value = '  -0.41%  '
$ws.Range("D18").Value = '2.968.87'
$ws.Range("E18").Value = '  -0.80%  '
$ws.Range("D19").Value = '59.053.45'
$ws.Range("E19").Value = '  +1.40%  '
$ws.Range("D20").Value = "'" + '434.48'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.61%  '
$ws.Range("D21").Value = "'" + '13.54'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.67%  '
$ws.Range("E22").Value = '  +0.75%  '
$ws.Range("D23").Value = "'" + '6.99'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.70%  '
$ws.Range("D24").Value = "'" + '13.04'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -3.36%  '
$ws.Range("D25").Value = "'" + '79.61'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.12%  '
$ws.Range("D26").Value = "'" + '1.00'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.08%  '
$ws.Range("E27").Value = '  +0.06%  '
$ws.Range("E28").Value = '  +6.50%  '
$ws.Range("E29").Value = '  +0.43%  '
$ws.Range("E30").Value = '  +0.59%  '
$ws.Range("D31").Value = "'" + '25.58'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.01%  '
$ws.Range("D32").Value = "'" + '6.15'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.62%  '
$ws.Range("E33").Value = '  +5.71%  '
$ws.Range("B34").Value = 'Mantle'
$ws.Range("C34").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D34").Value = "'" + '0.985'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.18%  '
$ws.Range("B35").Value = 'PEPE'
$ws.Range("C35").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D35").Value = '0.0₃0761'
$ws.Range("E35").Value = '  +5.55%  '
$ws.Range("B36").Value = 'Filecoin'
$ws.Range("C36").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D36").Value = "'" + '5.86'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.52%  '
$ws.Range("E37").Value = '  -2.27%  '
$ws.Range("D38").Value = "'" + '48.43'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.13%  '
$ws.Range("D39").Value = "'" + '8.68'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.58%  '
$ws.Range("D40").Value = "'" + '2.75'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.13%  '
$ws.Range("D41").Value = "'" + '394.79'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.87%  '
$ws.Range("D42").Value = "'" + '0.0349'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.40%  '
$ws.Range("D43").Value = '2.716.43'
$ws.Range("E43").Value = '  -0.66%  '
$ws.Range("D44").Value = "'" + '0.105'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.86%  '
$ws.Range("E45").Value = '  +2.02%  '
$ws.Range("E46").Value = '  -0.04%  '
$ws.Range("D47").Value = "'" + '122.26'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.38%  '
$ws.Range("D48").Value = "'" + '34.13'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +13.70%  '
$ws.Range("E49").Value = '  +0.13%  '
$ws.Range("E50").Value = '  -2.75%  '
$ws.Range("E51").Value = '  -1.21%  '
